# agregado de etiqueta 112 COGOTE BLOKE 10B JABAT - FALTA TRADUCCION
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New label data in row 24: etiqueta (label) number 112, translation pending ("FALTA")
$ws.Range("M24").Value = 112
$ws.Range("N24").Value = "FALTA"

# Widen column H to fit the new content (no longer auto "best fit")
$ws.Columns.Item(8).ColumnWidth = 26

# Restore the view: scroll so row 13 is at the top of the frozen pane,
# and leave the active selection on O24 (matches the saved workbook view)
$ws.Range("A13").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("O24").Select() | Out-Null
